$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yyymmdd")

# Update header cell G1 from "上繳日" to "繳息迄日"
$ws.Range("G1").Value = "繳息迄日"

# Reflect the updated selection/active cell as captured in the diff
$ws.Range("G1").Select()
